# Update cryptos list (price & 1h volume change) as scraped by GitHub Actions.
# Also fixes the ordering of two pairs of rows (Filecoin/FirstDigitalUSD and
# InjectiveProtocol/OKB swapped places in the ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.660.92"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.671.71"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "631.36"
$ws.Range("E5").Value = "  -6.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.59"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.497"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000230"
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.283.19"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.56"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.651.26"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.662.90"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.54"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.90"
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.30"
$ws.Range("E20").Value = "  +5.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "470.96"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.651"
$ws.Range("E22").Value = "  -0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.79"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.813.15"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.11"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.73"
$ws.Range("E28").Value = "  -4.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.60"
$ws.Range("E29").Value = "  -3.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.68"
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.165"
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.65"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.42"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.669.51"
$ws.Range("E36").Value = "  -0.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.34"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "178.72"
$ws.Range("E39").Value = "  +3.62%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.84"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0893"
$ws.Range("E43").Value = "  -1.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.926"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "29.14"
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.68"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.74"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.86"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000265"
$ws.Range("E49").Value = "  -4.86%  "
$ws.Range("E50").Value = "  -5.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.22"
$ws.Range("E51").Value = "  -5.47%  "
